# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1376.4445
$ws.Range("I53").Value = 2174.4
$ws.Range("J53").Value = 379
$ws.Range("K53").Value = 2174.4
$ws.Range("L53").Value = 379
$ws.Range("M53").Value = -1537.4
$ws.Range("N53").Value = -1653
$ws.Range("H70").Value = 3751
$ws.Range("J70").Value = 2500
$ws.Range("L70").Value = 7500
$ws.Range("N70").Value = -8040
$ws.Range("H73").Value = 3751
$ws.Range("J73").Value = 2500
$ws.Range("L73").Value = 7500
$ws.Range("N73").Value = -9372
$ws.Range("H80").Value = 1648.4445
$ws.Range("I80").Value = 1097.4
$ws.Range("J80").Value = 2337.25
$ws.Range("K80").Value = 3292.2
$ws.Range("L80").Value = 7011.75
$ws.Range("M80").Value = -2294.2
$ws.Range("N80").Value = -9007.75
$ws.Range("H82").Value = 1422.5
$ws.Range("I82").Value = 1422.5
$ws.Range("K82").Value = 4267.5
$ws.Range("M82").Value = -3861.5
$ws.Range("H83").Value = 1648.4445
$ws.Range("I83").Value = 1097.4
$ws.Range("J83").Value = 2337.25
$ws.Range("K83").Value = 9876.6
$ws.Range("L83").Value = 21035.25
$ws.Range("M83").Value = -4884.6
$ws.Range("N83").Value = -31019.25
$ws.Range("H85").Value = 1422.5
$ws.Range("I85").Value = 1422.5
$ws.Range("K85").Value = 4267.5
$ws.Range("M85").Value = -2863.5
$ws.Range("H107").Value = 5718.625
$ws.Range("I107").Value = 4291.5
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 4291.5
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = -2371.5
$ws.Range("N107").Value = -13840
$ws.Range("H111").Value = 1242.6666
$ws.Range("I111").Value = 1242.6666
$ws.Range("K111").Value = 3727.9998
$ws.Range("M111").Value = -660.9998000000001
$ws.Range("H113").Value = 8582.583000000001
$ws.Range("J113").Value = 8936.5
$ws.Range("L113").Value = 8936.5
$ws.Range("N113").Value = -15444.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 739.375
$ws.Range("I2").Value = 779.4545000000001
$ws.Range("J2").Value = 651.2
$ws.Range("K2").Value = 779.4545000000001
$ws.Range("L2").Value = 651.2
$ws.Range("M2").Value = -666.4545000000001
$ws.Range("N2").Value = -877.2
$ws.Range("H31").Value = 22499.75
$ws.Range("I31").Value = 22499.75
$ws.Range("K31").Value = 22499.75
$ws.Range("M31").Value = -22205.75
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H63").Value = 1269.2222
$ws.Range("I63").Value = 1302.875
$ws.Range("J63").Value = 1000
$ws.Range("K63").Value = 1302.875
$ws.Range("L63").Value = 1000
$ws.Range("M63").Value = -616.875
$ws.Range("N63").Value = -2372
$ws.Range("H66").Value = 1269.2222
$ws.Range("I66").Value = 1302.875
$ws.Range("J66").Value = 1000
$ws.Range("K66").Value = 6514.375
$ws.Range("L66").Value = 5000
$ws.Range("M66").Value = -3082.375
$ws.Range("N66").Value = -11864
$ws.Range("H74").Value = 4374
$ws.Range("I74").Value = 4332
$ws.Range("K74").Value = 4332
$ws.Range("M74").Value = -3458
$ws.Range("H77").Value = 4374
$ws.Range("I77").Value = 4332
$ws.Range("K77").Value = 21660
$ws.Range("M77").Value = -17292
$ws.Range("H116").Value = 739.375
$ws.Range("I116").Value = 779.4545000000001
$ws.Range("J116").Value = 651.2
$ws.Range("K116").Value = 779.4545000000001
$ws.Range("L116").Value = 651.2
$ws.Range("M116").Value = 1514.5455
$ws.Range("N116").Value = -5239.2
$ws.Range("H122").Value = 3199
$ws.Range("I122").Value = 2998.75
$ws.Range("K122").Value = 8996.25
$ws.Range("M122").Value = -6546.25
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 739.375
$ws.Range("I3").Value = 779.4545000000001
$ws.Range("J3").Value = 651.2
$ws.Range("K3").Value = 779.4545000000001
$ws.Range("L3").Value = 651.2
$ws.Range("M3").Value = -665.4545000000001
$ws.Range("N3").Value = -879.2
$ws.Range("H15").Value = 35964
$ws.Range("J15").Value = 35964
$ws.Range("L15").Value = 35964
$ws.Range("N15").Value = -36418
$ws.Range("H35").Value = 15000
$ws.Range("J35").Value = 15000
$ws.Range("L35").Value = 15000
$ws.Range("N35").Value = -15620
$ws.Range("H82").Value = 20161.143
$ws.Range("H85").Value = 20161.143
$ws.Range("H86").Value = 7451.5
$ws.Range("I86").Value = 7503
$ws.Range("K86").Value = 7503
$ws.Range("M86").Value = -6380
$ws.Range("H89").Value = 7451.5
$ws.Range("I89").Value = 7503
$ws.Range("K89").Value = 37515
$ws.Range("M89").Value = -31899
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 20541.5
$ws.Range("H59").Value = 28626
$ws.Range("I59").Value = 24504
$ws.Range("K59").Value = 24504
$ws.Range("M59").Value = -23359
$ws.Range("H60").Value = 22418.2
$ws.Range("I60").Value = 24030.334
$ws.Range("K60").Value = 24030.334
$ws.Range("M60").Value = -23519.334
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 49349.477
$ws.Range("I4").Value = 93253.25
$ws.Range("J4").Value = 1454.4546
$ws.Range("K4").Value = 279759.75
$ws.Range("L4").Value = 4363.3638
$ws.Range("M4").Value = -279647.75
$ws.Range("N4").Value = -4587.3638
$ws.Range("H39").Value = 2299.875
$ws.Range("J39").Value = 2299.875
$ws.Range("L39").Value = 6899.625
$ws.Range("N39").Value = -7487.625
$ws.Range("H55").Value = 366.2857
$ws.Range("I55").Value = 366.2857
$ws.Range("K55").Value = 1098.8571
$ws.Range("M55").Value = -921.8571000000002
$ws.Range("H107").Value = 947
$ws.Range("I107").Value = 895
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 2685
$ws.Range("L107").Value = 2997
$ws.Range("M107").Value = -765
$ws.Range("N107").Value = -6837
$ws.Range("H110").Value = 21451
$ws.Range("I110").Value = 21451
$ws.Range("K110").Value = 64353
$ws.Range("M110").Value = -60263
$ws.Range("H132").Value = 2250
$ws.Range("J132").Value = 2500
$ws.Range("L132").Value = 22500
$ws.Range("N132").Value = -27560
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 188.5
$ws.Range("I2").Value = 279.66666
$ws.Range("J2").Value = 71.28570999999999
$ws.Range("K2").Value = 279.66666
$ws.Range("L2").Value = 71.28570999999999
$ws.Range("M2").Value = -166.66666
$ws.Range("N2").Value = -297.28571
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3750
$ws.Range("I132").Value = 3750
$ws.Range("K132").Value = 11250
$ws.Range("M132").Value = -8720
$ws.Range("H136").Value = 3731.6
$ws.Range("I136").Value = 3177.75
$ws.Range("J136").Value = 5947
$ws.Range("K136").Value = 9533.25
$ws.Range("L136").Value = 17841
$ws.Range("M136").Value = -6983.25
$ws.Range("N136").Value = -22941
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 74998
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H81").Value = 999.5
$ws.Range("I81").Value = 999.5
$ws.Range("K81").Value = 1999
$ws.Range("M81").Value = -938
$ws.Range("H84").Value = 999.5
$ws.Range("I84").Value = 999.5
$ws.Range("K84").Value = 9995
$ws.Range("M84").Value = -4691
$ws.Range("H136").Value = 9279.6
$ws.Range("I136").Value = 8900
$ws.Range("K136").Value = 26700
$ws.Range("M136").Value = -24150
